$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay text-typed (mirrors source inlineStr cells,
# avoiding Excel auto-converting numeric-looking strings like "231.12" into numbers).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.509.80'
$ws.Range('E2').Value = '  -1.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.041.31'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.12'
$ws.Range('E5').Value = '  -11.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.599'
$ws.Range('E6').Value = '  -1.16%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.17'
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.370'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.96'
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0746'
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.338.23'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.33'
$ws.Range('E14').Value = '  +1.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.96'
$ws.Range('E15').Value = '  -9.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.757'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.09'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.023.36'
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '36.626.20'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.82'
$ws.Range('E20').Value = '  +15.27%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '67.50'
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0794'
$ws.Range('E22').Value = '  -3.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '220.25'
$ws.Range('E23').Value = '  -5.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E25').Value = '  +1.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.36'
$ws.Range('E26').Value = '  -8.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.78'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.69'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.126'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '18.90'
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.34'
$ws.Range('E31').Value = '  +3.50%  '
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.36'
$ws.Range('E33').Value = '  -3.08%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0601'
$ws.Range('E34').Value = '  -2.62%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.47'
$ws.Range('E35').Value = '  +4.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.25'
$ws.Range('E36').Value = '  -1.55%  '
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.75'
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.80'
$ws.Range('E39').Value = '  +8.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.23'
$ws.Range('E40').Value = '  -5.48%  '
$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.55'
$ws.Range('E41').Value = '  +46.34%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.95'
$ws.Range('E42').Value = '  -3.96%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.479.50'
$ws.Range('E43').Value = '  +3.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0938'
$ws.Range('E44').Value = '  +3.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.11'
$ws.Range('E45').Value = '  +5.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0203'
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.11'
$ws.Range('E47').Value = '  -4.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.54'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.89'
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.92'
$ws.Range('E51').Value = '  +2.01%  '
